$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# --- Fix typo in Q7 header: "Battery Sytandby" -> "Battery Standby" ---
$ws.Range("Q7").Value = "Battery Standby"

# --- Add new "Loading Detail" header/value columns at F1:G2 (previously blank) ---
# Move the values that used to live in W7:X8 up into F1:G2, matching their
# existing formatting (header style from row 7, data style from row 8).
$ws.Range("F1").Value = $ws.Range("W7").Value2
$ws.Range("G1").Value = $ws.Range("X7").Value2
$ws.Range("F2").Value = $ws.Range("W8").Value2
$ws.Range("G2").Value = $ws.Range("X8").Value2

$ws.Range("A7").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B8").Copy() | Out-Null
$ws.Range("F2:G2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Remove the now-obsolete W7:X8 columns of data ---
$ws.Range("W7:X8").Clear()

# --- Widen column G to fit the new content ---
$ws.Columns.Item(7).ColumnWidth = 18.8

# --- Update the saved view/selection for the sheet ---
$ws.Range("F1:G2").Select()

$wb.Save()
